$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'This file records the analysis results that I get when analysing data as shown in the Jupyter Notebook at: "C:\Users\TamCoding\Documents\Team_project\shared_git_repo\team-chip-project\data_analysis" '
$ws.Range("A2").Value = "The analysis aims to take wavelength and power transmission scan data for different structures, to find resonance peaks in this data and then calculate the free spectral range of the structure from this."
$ws.Range("A3").Value = "This is done by finding peaks using scipy.signal.find_peaks(), and then finding the difference in wavelength between these peak locations. This gives a series of wavelength fsrs."
$ws.Range("A4").Value = "Then, this series of wavelength fsrs is used to find its mean. Then the standard error is found by using np.std with ddof=1 to get the standard deviation, and then dividing by square root of number of fsrs in the series."
$ws.Range("A5").Value = "We only use a section of the full data in the CSV files for the peak finding. This is detailed by start array index and end array index (and also start wavelength and end wavelength)."
$ws.Range("A6").Value = "The files under team-chip-project are in a Git repo. At the moment working on branch: Tamanna19Dec_analysis."

$ws.Range("A8").Value = "Definitions: Data CSV filename: the filename of the scan data read in when doing the analysis.  "
$ws.Range("A9").Value = "Wavelength step size: the step size passed into the give_peak_locs function used to calculate the approx_fsr in number of data points rather than wavelength to be passed into scipy.signal.find_peaks() to be used as wlen."
$ws.Range("A10").Value = "Start array index: The index that the arrays of powers and wavelengths used to find resonance peaks starts at ( array includes this item). Note index is for numpy arrayso starts at 0"
$ws.Range("A11").Value = "End array index: The index that arrays of powers and wavelengths used to find resonance peaks ends one before (array does not include this index item)."
$ws.Range("A12").Value = "Start wavelength: The wavelength that the data used to find peaks starts at. This data point is included in data to find resonance peaks."
$ws.Range("A13").Value = "End wavelength: The wavelength that the data used to find peaks ends at. This data point is included in the data to find resonance peaks."
$ws.Range("A14").Value = "fsr_mean: the mean free spectral range as calculated by the code from finding resonance peaks as detailed above."
$ws.Range("A15").Value = "fsr_std_error: the standard error in the free spectral range, as calculated by the code from finding resonance peaks as detailed above."
$ws.Range("A16").Value = "prominence: the value passed into scipy.signal.find_peaks() when peak finding as the prominence parameter. See notes for exact definition of this. Chosen by looking by eye and picking about half of the height span of a noise part of the graph that doesn't contain the resonance peaks."
$ws.Range("A17").Value = "distance:the value passed into scipy.signal.find_peaks when peak finding as the distance parameter. The minimum distance between peaks in number of points. Choose this by taking approx_fsr/2 and dividing that by wavelength step size."
$ws.Range("A18").Value = "approx_fsr: chosen by looking at plot, and maybe using a ruler, and guessing by eye roughly what the free spectral range is. It is divided by the wavelength step size to give the wlen parameter to be passed into scipy.signal.find_peaks() for finding the peaks. "
$ws.Range("A19").Value = "The procedure will be to Git commit after each run and after I've added the parameters to this table, so each commit should contain the state of the code and state of the table after each run has been recorded."

$ws.Range("A21").Value = "Data CSV filename"
$ws.Range("B21").Value = "Wavelength step size"
$ws.Range("C21").Value = "Start array index"
$ws.Range("D21").Value = "End array index"
$ws.Range("E21").Value = "Start wavelength"
$ws.Range("F21").Value = "End wavelength"
$ws.Range("G21").Value = "prominence"
$ws.Range("H21").Value = "distance"
$ws.Range("I21").Value = "approx_fsr"
$ws.Range("J21").Value = "fsr_mean"
$ws.Range("K21").Value = "fsr_std error"

$ws.Range("B21:F21").EntireColumn.AutoFit() | Out-Null
$ws.Range("G21:I21").ColumnWidth = 14.08984375
$ws.Range("K21:L21").EntireColumn.AutoFit() | Out-Null

$ws.Range("A15").Select() | Out-Null
